{"js": "// Update the worksheet date and every three-digit x one-digit multiplication\n// problem/answer pair with the newly generated values.\nconst replacements = [\n  [\"2026-01-30 Friday\", \"2026-01-31 Saturday\"],\n  [\"361\u00d75=1805\", \"590\u00d72=1180\"],\n  [\"996\u00d78=7968\", \"569\u00d78=4552\"],\n  [\"937\u00d76=5622\", \"413\u00d77=2891\"],\n  [\"808\u00d79=7272\", \"793\u00d74=3172\"],\n  [\"865\u00d76=5190\", \"712\u00d72=1424\"],\n  [\"353\u00d76=2118\", \"282\u00d73=846\"],\n  [\"950\u00d74=3800\", \"737\u00d79=6633\"],\n  [\"644\u00d76=3864\", \"106\u00d74=424\"],\n  [\"159\u00d72=318\", \"119\u00d73=357\"],\n  [\"266\u00d74=1064\", \"508\u00d79=4572\"],\n  [\"972\u00d75=4860\", \"845\u00d76=5070\"],\n  [\"464\u00d75=2320\", \"913\u00d72=1826\"],\n  [\"761\u00d77=5327\", \"281\u00d75=1405\"],\n  [\"410\u00d72=820\", \"592\u00d79=5328\"],\n  [\"448\u00d72=896\", \"548\u00d77=3836\"],\n  [\"519\u00d78=4152\", \"966\u00d76=5796\"],\n  [\"857\u00d75=4285\", \"270\u00d79=2430\"],\n  [\"109\u00d75=545\", \"196\u00d77=1372\"],\n  [\"968\u00d73=2904\", \"452\u00d76=2712\"],\n  [\"163\u00d76=978\", \"182\u00d73=546\"],\n  [\"740\u00d78=5920\", \"520\u00d79=4680\"],\n  [\"409\u00d76=2454\", \"437\u00d78=3496\"],\n  [\"277\u00d75=1385\", \"213\u00d77=1491\"],\n  [\"202\u00d74=808\", \"652\u00d73=1956\"],\n  [\"664\u00d79=5976\", \"859\u00d73=2577\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true });\n  results.load(\"text\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(`Text not found: ${oldText}`);\n  }\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the worksheet date and every three-digit x one-digit multiplication\n# problem/answer pair with the newly generated values.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2026-01-30 Friday\", \"2026-01-31 Saturday\"),\n    @(\"361\u00d75=1805\", \"590\u00d72=1180\"),\n    @(\"996\u00d78=7968\", \"569\u00d78=4552\"),\n    @(\"937\u00d76=5622\", \"413\u00d77=2891\"),\n    @(\"808\u00d79=7272\", \"793\u00d74=3172\"),\n    @(\"865\u00d76=5190\", \"712\u00d72=1424\"),\n    @(\"353\u00d76=2118\", \"282\u00d73=846\"),\n    @(\"950\u00d74=3800\", \"737\u00d79=6633\"),\n    @(\"644\u00d76=3864\", \"106\u00d74=424\"),\n    @(\"159\u00d72=318\", \"119\u00d73=357\"),\n    @(\"266\u00d74=1064\", \"508\u00d79=4572\"),\n    @(\"972\u00d75=4860\", \"845\u00d76=5070\"),\n    @(\"464\u00d75=2320\", \"913\u00d72=1826\"),\n    @(\"761\u00d77=5327\", \"281\u00d75=1405\"),\n    @(\"410\u00d72=820\", \"592\u00d79=5328\"),\n    @(\"448\u00d72=896\", \"548\u00d77=3836\"),\n    @(\"519\u00d78=4152\", \"966\u00d76=5796\"),\n    @(\"857\u00d75=4285\", \"270\u00d79=2430\"),\n    @(\"109\u00d75=545\", \"196\u00d77=1372\"),\n    @(\"968\u00d73=2904\", \"452\u00d76=2712\"),\n    @(\"163\u00d76=978\", \"182\u00d73=546\"),\n    @(\"740\u00d78=5920\", \"520\u00d79=4680\"),\n    @(\"409\u00d76=2454\", \"437\u00d78=3496\"),\n    @(\"277\u00d75=1385\", \"213\u00d77=1491\"),\n    @(\"202\u00d74=808\", \"652\u00d73=1956\"),\n    @(\"664\u00d79=5976\", \"859\u00d73=2577\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $newText\n    $find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, $newText, 2)\n}\n"}
